$wb = $excel.ActiveWorkbook

# --- Rename the "Include #0" sheet ---
$wsMeta = $wb.Worksheets.Item(1)
$wsInc  = $wb.Worksheets.Item(2)
$wsInc.Name = "Include from International Cl"

# --- Sheet1 ("Metadata") updates ---
# Rows 1-9 are unchanged. Rows 10 onward change values and row 16 is removed
# (the sheet shrinks from 16 to 15 rows).

$wsMeta.Range("B3").Value  = "0.1.6"
$wsMeta.Range("B6").Value  = "active"
$wsMeta.Range("B8").Value  = "2023-05-05T10:50:04-05:00"
$wsMeta.Range("B10").Value = "No display for ContactDetail"
$wsMeta.Range("B11").Value = "No display for ContactDetail"
$wsMeta.Range("A12").Value = "Description"
$wsMeta.Range("B12").Value = "GVHD ICD-10 Codes"
$wsMeta.Range("A13").Value = "Purpose"
$wsMeta.Range("B13").Value = ""
$wsMeta.Range("A14").Value = "Copyright"
$wsMeta.Range("B14").Value = ""
$wsMeta.Range("A15").Value = "Immutable"
$wsMeta.Range("B15").Value = "BooleanType[null]"

# Drop the now-superfluous last row (previously row 16: Operation-related rows
# belong to the other sheet; this row is simply removed by the revert).
$wsMeta.Rows.Item(16).Delete()

# --- Sheet2 ("Include from International Cl") updates ---
# Row 3 (A3/B3) is already an empty string in both the before and after
# states, so it is left untouched here.
$wsInc.Range("B1").Value = "Operation"
$wsInc.Range("A2").Value = "concept"
$wsInc.Range("B2").Value = "is-a"
$wsInc.Range("C2").Value = "D89.81"
$wsInc.Range("A4").Value = "System URI"
$wsInc.Range("B4").Value = "http://hl7.org/fhir/sid/icd-10-cm"
